$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 684.8571
$ws.Range("I28").Value = 568.625
$ws.Range("K28").Value = 568.625
$ws.Range("M28").Value = -83.625
$ws.Range("H47").Value = 30000
$ws.Range("I47").Value = 30000
$ws.Range("K47").Value = 30000
$ws.Range("M47").Value = -29028
$ws.Range("H75").Value = 75550
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 75550
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 75550
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -77422
$ws.Range("H78").Value = 75550
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 75550
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 226650
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -236010
$ws.Range("H98").Value = 1313.8572
$ws.Range("I98").Value = 1200.8889
$ws.Range("K98").Value = 1200.8889
$ws.Range("M98").Value = 297.1111000000001
$ws.Range("H107").Value = 1313.7
$ws.Range("I107").Value = 1624
$ws.Range("K107").Value = 1624
$ws.Range("M107").Value = 296
$ws.Range("H122").Value = 1313.8572
$ws.Range("I122").Value = 1200.8889
$ws.Range("K122").Value = 3602.6667
$ws.Range("M122").Value = -1152.6667
$ws.Range("H135").Value = 596.85
$ws.Range("I135").Value = 609.5897
$ws.Range("K135").Value = 5486.3073
$ws.Range("M135").Value = -2951.3073
$ws.Range("H136").Value = 119449.8
$ws.Range("J136").Value = 119449.8
$ws.Range("L136").Value = 119449.8
$ws.Range("N136").Value = -129649.8
$ws.Range("H137").Value = 2895.5833
$ws.Range("I137").Value = 2851.5173
$ws.Range("K137").Value = 8554.5519
$ws.Range("M137").Value = -6004.5519

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 17858406
$ws.Range("I2").Value = 19231736
$ws.Range("K2").Value = 19231736
$ws.Range("M2").Value = -19231623
$ws.Range("H28").Value = 12311.667
$ws.Range("J28").Value = 20000
$ws.Range("L28").Value = 20000
$ws.Range("N28").Value = -20384
$ws.Range("H43").Value = 25098.682
$ws.Range("J43").Value = 22403.125
$ws.Range("L43").Value = 22403.125
$ws.Range("N43").Value = -23029.125
$ws.Range("H74").Value = 3536.84
$ws.Range("I74").Value = 2045.9445
$ws.Range("K74").Value = 2045.9445
$ws.Range("M74").Value = -1171.9445
$ws.Range("H77").Value = 3536.84
$ws.Range("I77").Value = 2045.9445
$ws.Range("K77").Value = 10229.7225
$ws.Range("M77").Value = -5861.7225
$ws.Range("H86").Value = 21268
$ws.Range("J86").Value = 21268
$ws.Range("L86").Value = 21268
$ws.Range("N86").Value = -23640
$ws.Range("H89").Value = 21268
$ws.Range("J89").Value = 21268
$ws.Range("L89").Value = 63804
$ws.Range("N89").Value = -75660
$ws.Range("H99").Value = 12311.667
$ws.Range("J99").Value = 20000
$ws.Range("L99").Value = 20000
$ws.Range("N99").Value = -25990
$ws.Range("H116").Value = 17858406
$ws.Range("I116").Value = 19231736
$ws.Range("K116").Value = 19231736
$ws.Range("M116").Value = -19229442
$ws.Range("H132").Value = 1197.3684
$ws.Range("I132").Value = 1093.7059
$ws.Range("K132").Value = 3281.1177
$ws.Range("M132").Value = -751.1176999999998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 17858406
$ws.Range("I3").Value = 19231736
$ws.Range("K3").Value = 19231736
$ws.Range("M3").Value = -19231622
$ws.Range("H86").Value = 1804.8334
$ws.Range("I86").Value = 1449.8334
$ws.Range("J86").Value = 2159.8333
$ws.Range("K86").Value = 1449.8334
$ws.Range("L86").Value = 2159.8333
$ws.Range("M86").Value = -326.8334
$ws.Range("N86").Value = -4405.8333
$ws.Range("H89").Value = 1804.8334
$ws.Range("I89").Value = 1449.8334
$ws.Range("J89").Value = 2159.8333
$ws.Range("K89").Value = 7249.166999999999
$ws.Range("L89").Value = 10799.1665
$ws.Range("M89").Value = -1633.166999999999
$ws.Range("N89").Value = -22031.1665
$ws.Range("H134").Value = 2289.3125
$ws.Range("I134").Value = 2289.3125
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6867.9375
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -4332.9375
$ws.Range("N134").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 78427
$ws.Range("J9").Value = 78427
$ws.Range("L9").Value = 78427
$ws.Range("N9").Value = -78763
$ws.Range("H99").Value = 5671.222
$ws.Range("I99").Value = 5935.2
$ws.Range("K99").Value = 5935.2
$ws.Range("M99").Value = -4437.2
$ws.Range("H104").Value = 32587.125
$ws.Range("J104").Value = 32587.125
$ws.Range("L104").Value = 32587.125
$ws.Range("N104").Value = -37829.125
$ws.Range("H105").Value = 1019.5
$ws.Range("I105").Value = 1051.25
$ws.Range("J105").Value = 924.25
$ws.Range("K105").Value = 1051.25
$ws.Range("L105").Value = 924.25
$ws.Range("M105").Value = 695.75
$ws.Range("N105").Value = -4418.25
$ws.Range("H126").Value = 5671.222
$ws.Range("I126").Value = 5935.2
$ws.Range("K126").Value = 17805.6
$ws.Range("M126").Value = -15335.6
$ws.Range("H132").Value = 3791.8462
$ws.Range("I132").Value = 2643.889
$ws.Range("K132").Value = 7931.667
$ws.Range("M132").Value = -5401.667

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 160.35715
$ws.Range("I61").Value = 77.85714
$ws.Range("K61").Value = 233.57142
$ws.Range("M61").Value = -18.57141999999999
$ws.Range("H101").Value = 11999
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H113").Value = 1926.9445
$ws.Range("I113").Value = 544.5
$ws.Range("J113").Value = 2099.75
$ws.Range("K113").Value = 1633.5
$ws.Range("L113").Value = 6299.25
$ws.Range("M113").Value = 536.5
$ws.Range("N113").Value = -10639.25
$ws.Range("H122").Value = 6695.6665
$ws.Range("I122").Value = 4033.1667
$ws.Range("J122").Value = 8470.667
$ws.Range("K122").Value = 36298.5003
$ws.Range("L122").Value = 76236.003
$ws.Range("M122").Value = -33848.5003
$ws.Range("N122").Value = -81136.003
$ws.Range("H129").Value = 464473.47
$ws.Range("I129").Value = 78028.08
$ws.Range("J129").Value = 921181.6
$ws.Range("K129").Value = 234084.24
$ws.Range("L129").Value = 2763544.8
$ws.Range("M129").Value = -229084.24
$ws.Range("N129").Value = -2773544.8
$ws.Range("H137").Value = 2064.2632
$ws.Range("I137").Value = 1753.7273
$ws.Range("J137").Value = 2491.25
$ws.Range("K137").Value = 5261.1819
$ws.Range("L137").Value = 7473.75
$ws.Range("M137").Value = -161.1818999999996
$ws.Range("N137").Value = -17673.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 33.25926
$ws.Range("I2").Value = 35.708332
$ws.Range("K2").Value = 35.708332
$ws.Range("M2").Value = 77.291668
$ws.Range("H132").Value = 1692.775
$ws.Range("I132").Value = 897.21875
$ws.Range("K132").Value = 2691.65625
$ws.Range("M132").Value = -161.65625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 18773.615
$ws.Range("I93").Value = 3594.75
$ws.Range("J93").Value = 43059.8
$ws.Range("K93").Value = 3594.75
$ws.Range("L93").Value = 43059.8
$ws.Range("M93").Value = -2346.75
$ws.Range("N93").Value = -45555.8
$ws.Range("H98").Value = 65000
$ws.Range("J98").Value = 65000
$ws.Range("L98").Value = 65000
$ws.Range("N98").Value = -70990
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H139").Value = 67356.336

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3237.125
$ws.Range("J62").Value = 2985.2856
$ws.Range("L62").Value = 2985.2856
$ws.Range("N62").Value = -4233.2856
$ws.Range("H65").Value = 3237.125
$ws.Range("J65").Value = 2985.2856
$ws.Range("L65").Value = 14926.428
$ws.Range("N65").Value = -21166.428
$ws.Range("H96").Value = 2750.762
$ws.Range("J96").Value = 2652.1538
$ws.Range("L96").Value = 2652.1538
$ws.Range("N96").Value = -5398.1538
$ws.Range("H100").Value = 1115.0834
$ws.Range("I100").Value = 1002.8571
$ws.Range("K100").Value = 2005.7142
$ws.Range("M100").Value = -1464.7142
$ws.Range("H113").Value = 1102.3684
$ws.Range("I113").Value = 836.80646
$ws.Range("J113").Value = 2278.4285
$ws.Range("K113").Value = 2510.41938
$ws.Range("L113").Value = 6835.2855
$ws.Range("M113").Value = -340.4193800000003
$ws.Range("N113").Value = -11175.2855
